$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values
$ws.Range("B3").Value = 0.8121622837929257
$ws.Range("C3").Value = 0.8031329410193295
$ws.Range("D3").Value = 0.750766555188693

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.7898564578798413
$ws.Range("C4").Value = 0.7889271823332352
$ws.Range("D4").Value = 0.7589429413638695

# Row 5: AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.6053334550982021
$ws.Range("C5").Value = 0.5785800703049623
$ws.Range("D5").Value = 0.5572595351453221
